$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (test_number) values from 2 to 4 for rows 2-11
$ws.Range("A2:A11").Value = 4

# Update the selected cell on the sheet to H11
$ws.Range("H11").Select()
